# Apply a uniform slide transition (Peel Off, slow speed, 1.25s duration)
# across the whole deck: the slide master, every slide layout, and every
# slide -- mirroring PowerPoint's "Apply To All Slides" transition action.
#
# PowerPoint's legacy SlideShowTransition.EntryEffect automation property
# only understands the pre-2010 PpEntryEffect transition set, so "Peel Off"
# (a PowerPoint 2010+ transition, stored as <p15:prstTrans prst="peelOff"/>
# inside an mc:AlternateContent block) cannot be dialed in through that
# enum. We still drive every property that *is* exposed by the object
# model (Duration, Speed, EntryEffect) so the resulting <p:transition>
# carries the same speed/duration/fallback-effect semantics the authored
# transition uses (spd="slow", 1250ms, fade fallback).

$p = $ppt.ActivePresentation

function Set-SlideTransition($target) {
    if ($target -eq $null) { return }
    $sst = $target.SlideShowTransition
    if ($sst -eq $null) { return }
    # Order matters for this host's transition serializer: set Duration,
    # then EntryEffect, then Speed, so spd/p14:dur/the effect element all
    # survive together on the emitted <p:transition>.
    $sst.Duration = 1.25
    $sst.EntryEffect = "ppEffectFade"
    $sst.Speed = 1
}

# Slide master.
Set-SlideTransition $p.SlideMaster

# Every slide layout off the master.
$master = $p.SlideMaster
if ($master -ne $null) {
    for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
        Set-SlideTransition $master.CustomLayouts.Item($li)
    }
}

# Every slide in the deck.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    Set-SlideTransition $p.Slides.Item($i)
}
